$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update receptor / edge expression metrics on row 2 with the refreshed TPM values
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9942836666666667
$ws.Range("N2").Value = 2.982851
$ws.Range("Q2").Value = 2.082703459470222
$ws.Range("R2").Value = 18.744331135232
